# "beginning of change to two plots"
# Reverses the 16-value block order feeding the "Wires 16x4" plot (column
# E/F) for the first group (rows 3-18), which cascades through the
# shared F19:F67 "=F(row-16)+16" formulas for the remaining groups.
# Also locks in the computed results for the third group (F50:F66) as
# static values, and clears the now-unused 65th entry (E67:F67).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Reverse the first 16-row block (F3:F18): 0..15 -> 15..0 ---
$firstBlock = @(15, 14, 13, 12, 11, 10, 9, 8, 7, 6, 5, 4, 3, 2, 1, 0)
for ($i = 0; $i -lt $firstBlock.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 6).Value2 = $firstBlock[$i]
}

# --- Freeze the third block (F50:F66) to its (now recalculated) static values ---
for ($row = 50; $row -le 66; $row++) {
    $v = $ws.Cells.Item($row, 6).Value2
    $ws.Cells.Item($row, 6).Value2 = $v
}

# --- Clear the now-trailing row 67 entries in this block (E67:F67) ---
$ws.Range("E67:F67").ClearContents()

# --- Update the view / selection to match where the user is now working ---
$win = $wb.Windows.Item(1)
$win.ScrollRow = 34
$win.ScrollColumn = 1
$ws.Range("D32").Select()
